$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the "Espárragos" (Feria Lagunitas de
# Puerto Montt) series. It belongs right after the header/earlier rows, at
# what becomes row 63 - so insert a fresh row there, which pushes the old
# rows 63-71 down to 64-72 (and grows the used range to A1:R72).
$ws.Rows("63:63").Insert()

# Populate the newly inserted row 63 with the new observation's data.
$ws.Range("A63").Value = 4
$ws.Range("B63").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C63").Value = "Los Lagos"
$ws.Range("D63").Value = 45209
$ws.Range("E63").Value = 10
$ws.Range("F63").Value = 300000000
$ws.Range("G63").Value = "Espárragos"
$ws.Range("H63").Value = "Sin especificar"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 400
$ws.Range("K63").Value = 1800
$ws.Range("L63").Value = 2000
$ws.Range("M63").Value = 1900
$ws.Range("N63").Value = "`$/kilo"
$ws.Range("O63").Value = "Provincia de Linares"
$ws.Range("P63").Value = 1900
$ws.Range("Q63").Value = 1
$ws.Range("R63").Value = "Hortaliza"
